$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F3").Value = 1107
$ws1.Range("F4").Value = 1915
$ws1.Range("F5").Value = 596
$ws1.Range("F6").Value = 1247
$ws1.Range("F8").Value = 34
$ws1.Range("F10").Value = 330
$ws1.Range("F11").Value = 113
$ws1.Range("F13").Value = 806
$ws1.Range("F14").Value = 239
$ws1.Range("F15").Value = 124
$ws1.Range("F19").Value = 221
$ws1.Range("F20").Value = 694
$ws1.Range("F21").Value = 71
$ws1.Range("F22").Value = 662
$ws1.Range("F23").Value = 186
$ws1.Range("F25").Value = 904
$ws1.Range("F26").Value = 351
$ws1.Range("F29").Value = 303
$ws1.Range("F31").Value = 20
$ws1.Range("F32").Value = 422
$ws2.Range("F4").Value = 333
$ws2.Range("F7").Value = 264
$ws4.Range("F4").Value = 1107
$ws4.Range("F5").Value = 1915
$ws4.Range("F6").Value = 596
$ws4.Range("F7").Value = 1247
$ws4.Range("F10").Value = 34
$ws4.Range("F12").Value = 330
$ws4.Range("F13").Value = 113
$ws4.Range("F15").Value = 806
$ws4.Range("F16").Value = 239
$ws4.Range("F17").Value = 124
$ws4.Range("F20").Value = 333
$ws4.Range("F25").Value = 264
$ws4.Range("F26").Value = 264
$ws4.Range("F27").Value = 221
$ws4.Range("F28").Value = 694
$ws4.Range("F29").Value = 71
$ws4.Range("F30").Value = 662
$ws4.Range("F31").Value = 186
$ws4.Range("F33").Value = 904
$ws4.Range("F34").Value = 351
$ws4.Range("F39").Value = 303
$ws4.Range("F44").Value = 20
$ws4.Range("F46").Value = 422
